$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# Update the "Lower Right Cell" values in column D to reflect the new
# scenario listing extent (row 24 -> row 25) for the new scenario (29).
$ws.Range("D5").Value = "A25"
$ws.Range("D6").Value = "B25"
$ws.Range("D7").Value = "C25"
$ws.Range("D8").Value = "G25"
$ws.Range("D9").Value = "H25"
$ws.Range("D10").Value = "I25"
$ws.Range("D11").Value = "J25"

# Update the active selection to match the saved workbook state.
$ws.Range("D12").Select()
